$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coded year Pathogen")

# Capture the existing threaded comment's text before we shift columns.
$oldCommentRange = $ws.Range("J116")
$oldComment = $oldCommentRange.Comment
$commentText = $oldComment.Text()

# Delete column D (the "Cholera" column). Everything to the right shifts one column left.
$ws.Range("D:D").Delete()

# Rename the "Russian flu pandemic (H1N1)" header (now at E1) to "H1N1 Russian Flu".
$ws.Range("E1").Value = "H1N1 Russian Flu"

# The comment cell (old J116) is now I116 after the column shift; recreate it there.
$oldComment.Delete()
$ws.Range("I116").AddCommentThreaded($commentText) | Out-Null

Write-Host "D1:" $ws.Range("D1").Text
Write-Host "E1:" $ws.Range("E1").Text
Write-Host "I116 comment:" $ws.Range("I116").Comment.Text()
